$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 367.33334
$ws.Range("J6").Value = 367.33334
$ws.Range("L6").Value = 1102.00002
$ws.Range("N6").Value = -1326.00002

$ws.Range("H15").Value = 182.15
$ws.Range("I15").Value = 182.15
$ws.Range("K15").Value = 546.45
$ws.Range("M15").Value = -377.45

$ws.Range("H33").Value = 399.86667
$ws.Range("I33").Value = 358.6842
$ws.Range("K33").Value = 358.6842
$ws.Range("M33").Value = -129.6842

$ws.Range("H86").Value = 129500.375
$ws.Range("I86").Value = 253750.75
$ws.Range("J86").Value = 5250
$ws.Range("K86").Value = 253750.75
$ws.Range("L86").Value = 5250
$ws.Range("M86").Value = -252627.75
$ws.Range("N86").Value = -7496

$ws.Range("H89").Value = 129500.375
$ws.Range("I89").Value = 253750.75
$ws.Range("J89").Value = 5250
$ws.Range("K89").Value = 1268753.75
$ws.Range("L89").Value = 26250
$ws.Range("M89").Value = -1263137.75
$ws.Range("N89").Value = -37482

$ws.Range("H106").Value = 1196.5
$ws.Range("I106").Value = 1196.5
$ws.Range("K106").Value = 1196.5
$ws.Range("M106").Value = -565.5

$ws.Range("H138").Value = 3590.6914
$ws.Range("I138").Value = 2732.8823
$ws.Range("J138").Value = 3780.078
$ws.Range("K138").Value = 8198.6469
$ws.Range("L138").Value = 11340.234
$ws.Range("M138").Value = -3058.6469
$ws.Range("N138").Value = -21620.234

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 14801.091
$ws.Range("I3").Value = 10402.5
$ws.Range("J3").Value = 15778.556
$ws.Range("K3").Value = 10402.5
$ws.Range("L3").Value = 15778.556
$ws.Range("M3").Value = -10287.5
$ws.Range("N3").Value = -16008.556

$ws.Range("H32").Value = 7759.108
$ws.Range("I32").Value = 7244.5425
$ws.Range("J32").Value = 12819
$ws.Range("K32").Value = 7244.5425
$ws.Range("L32").Value = 12819
$ws.Range("M32").Value = -6957.5425
$ws.Range("N32").Value = -13393

$ws.Range("H63").Value = 10500.5
$ws.Range("I63").Value = 3001
$ws.Range("J63").Value = 18000
$ws.Range("K63").Value = 3001
$ws.Range("L63").Value = 18000
$ws.Range("M63").Value = -2315
$ws.Range("N63").Value = -19372

$ws.Range("H66").Value = 10500.5
$ws.Range("I66").Value = 3001
$ws.Range("J66").Value = 18000
$ws.Range("K66").Value = 15005
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -11573
$ws.Range("N66").Value = -96864

$ws.Range("H140").Value = 41260
$ws.Range("I140").Value = 30390
$ws.Range("K140").Value = 30390
$ws.Range("M140").Value = -25210

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 37000
$ws.Range("J35").Value = 37000
$ws.Range("L35").Value = 37000
$ws.Range("N35").Value = -37620

$ws.Range("H82").Value = 21267.5
$ws.Range("I82").Value = 7257
$ws.Range("J82").Value = 35278
$ws.Range("K82").Value = 7257
$ws.Range("L82").Value = 35278
$ws.Range("M82").Value = -6874
$ws.Range("N82").Value = -36044

$ws.Range("H85").Value = 21267.5
$ws.Range("I85").Value = 7257
$ws.Range("J85").Value = 35278
$ws.Range("K85").Value = 7257
$ws.Range("L85").Value = 35278
$ws.Range("M85").Value = -5931
$ws.Range("N85").Value = -37930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1874.3334
$ws.Range("I16").Value = 1883.625
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 1883.625
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -1596.625
$ws.Range("N16").Value = -2374

$ws.Range("H31").Value = 2098.6262
$ws.Range("I31").Value = 1446
$ws.Range("J31").Value = 3403.879
$ws.Range("K31").Value = 1446
$ws.Range("L31").Value = 3403.879
$ws.Range("M31").Value = -1151
$ws.Range("N31").Value = -3993.879

$ws.Range("H34").Value = 2098.6262
$ws.Range("I34").Value = 1446
$ws.Range("J34").Value = 3403.879
$ws.Range("K34").Value = 1446
$ws.Range("L34").Value = 3403.879
$ws.Range("M34").Value = -1244
$ws.Range("N34").Value = -3807.879

$ws.Range("H39").Value = 7019
$ws.Range("I39").Value = 3546.375
$ws.Range("K39").Value = 3546.375
$ws.Range("M39").Value = -3155.375

$ws.Range("H49").Value = 7019
$ws.Range("I49").Value = 3546.375
$ws.Range("K49").Value = 3546.375
$ws.Range("M49").Value = -3364.375

$ws.Range("H111").Value = 78266.336
$ws.Range("J111").Value = 78266.336
$ws.Range("L111").Value = 78266.336
$ws.Range("N111").Value = -86446.336

$ws.Range("H113").Value = 1874.3334
$ws.Range("I113").Value = 1883.625
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1883.625
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 286.375
$ws.Range("N113").Value = -6140

$ws.Range("H122").Value = 7960.091
$ws.Range("I122").Value = 3861.3333
$ws.Range("J122").Value = 11375.723
$ws.Range("K122").Value = 11583.9999
$ws.Range("L122").Value = 34127.169
$ws.Range("M122").Value = -9133.999899999999
$ws.Range("N122").Value = -39027.169

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 800
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2706.4167
$ws.Range("I102").Value = 2413.8635
$ws.Range("K102").Value = 2413.8635
$ws.Range("M102").Value = -791.8634999999999

$ws.Range("H122").Value = 5300.1875
$ws.Range("I122").Value = 7828.5713
$ws.Range("J122").Value = 3333.6667
$ws.Range("K122").Value = 23485.7139
$ws.Range("L122").Value = 10001.0001
$ws.Range("M122").Value = -21035.7139
$ws.Range("N122").Value = -14901.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7017.207
$ws.Range("I122").Value = 5388.1763
$ws.Range("J122").Value = 9325
$ws.Range("K122").Value = 16164.5289
$ws.Range("L122").Value = 27975
$ws.Range("M122").Value = -13714.5289
$ws.Range("N122").Value = -32875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3420.4443
$ws.Range("I122").Value = 2131.625
$ws.Range("J122").Value = 4451.5
$ws.Range("K122").Value = 6394.875
$ws.Range("L122").Value = 13354.5
$ws.Range("M122").Value = -3944.875
$ws.Range("N122").Value = -18254.5
